$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, pushing existing rows 33:55 down to 34:56.
# Excel copies the formatting of the row above by default for the new row.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly price-report entry.
# Non-numeric / categorical columns mirror the entry that used to sit in row 33
# (now shifted to row 34): Mercado ID, Mercado, Región, Codreg, Categoría ID,
# Categoría, Variedad, Calidad, Unidad de comercialización, Origen,
# Kg o Unidades and Clasificación stay the same.
$ws.Cells.Item(33, 1).Value = 11
$ws.Cells.Item(33, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(33, 3).Value = "Bíobío"
$ws.Cells.Item(33, 4).Value = "01/21/2022"
$ws.Cells.Item(33, 5).Value = 8
$ws.Cells.Item(33, 6).Value = 100112012
$ws.Cells.Item(33, 7).Value = "Espinaca"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 50
$ws.Cells.Item(33, 11).Value = 8000
$ws.Cells.Item(33, 12).Value = 8500
$ws.Cells.Item(33, 13).Value = 8200
$ws.Cells.Item(33, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(33, 15).Value = "Región Metropolitana"
$ws.Cells.Item(33, 16).Value = 820
$ws.Cells.Item(33, 17).Value = 10
$ws.Cells.Item(33, 18).Value = "Hortaliza"
